$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections on the keyboard-shortcuts help page ---

# Row 16 (A16): Windows shortcut for removing background changed to CTRL+h
$ws.Range("A16").Value = "h\u00A0\u00A0\u00A0[Windows: CTRL+h]"

# Row 17 (A17): Windows shortcut for remove background curve changed from CTRL+ALT+h to CTRL+SHIFT+h
$ws.Range("A17").Value = "ALT+h\u00A0\u00A0\u00A0[Windows: CTRL+SHIFT+h]"

# Row 23: quick special-event shortcut text + expanded description
$ws.Range("A23").Value = "q,w,e,r + <value>"
$ws.Range("B23").Value = "Quick Special Event Entry.  The keys q,w,e, and r correspond to special events 1,2,3 and 4.  A two digit numeric value must follow the shortcut letter, e.g. 'q75', when the correspoding event slider max value is 100 or less (default setting).   When the slider max value is greater than 100, three digits must be entered and for values less that 100 a leading zero is required, e.g. 'q075'.  "

# Row 24: quick PID SV shortcut text + expanded description
$ws.Range("A24").Value = "v +  <value>"
$ws.Range("B24").Value = "Quick PID SV Entry.  Value is a three digit number.  For values less than 100 must be entered with a leading zero, e.g. 'v075'."

# --- View state: scroll position and selection changed ---
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B24").Select()

# --- Page setup: printing orientation set to portrait ---
$ws.PageSetup.Orientation = 1
